# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" worksheet (cloned from "2022-Q3" so it keeps the
#    exact same look/formatting), placed right after "总计" and before
#    "2022-Q3", populated with the new quarter's numbers.
# 2) Insert a new top row into "总计" for the new quarter, pushing the
#    existing quarter rows down (their own values are untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q4" sheet, cloned from "2022-Q3"
# ---------------------------------------------------------------------
$zj  = $wb.Worksheets.Item("总计")
$src = $wb.Worksheets.Item("2022-Q3")

$src.Copy($null, $zj)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# D2:G3 hold numbers formatted as plain text in this workbook (t="inlineStr"
# in the original) - force text storage via NumberFormat, write the text,
# then drop the format/style again so the cells end up style-less like the
# source cells, only keeping their "this is text" nature.
$numericLookingCells = $q4.Range("D2:G3")
$numericLookingCells.NumberFormat = "@"

$q4.Range("D2:D3").Value = "1.79"
$q4.Range("E2:E3").Value = "88.58"
$q4.Range("F2:F3").Value = "4.17"
$q4.Range("G2:G3").Value = "0.0746"
$q4.Range("H2:H3").Value = 7

$numericLookingCells.Style = "Normal"

# ---------------------------------------------------------------------
# Step 2: new top data row in "总计" for 2022-Q4
# ---------------------------------------------------------------------
$zj.Rows.Item(2).Insert()

# Insert() drags the header row's bold/border formatting onto the new row;
# strip it back to plain (matches how the other data rows look).
$zj.Range("B2:D2").Style = "Normal"

# Column A carries the bold/bordered "index" style (s="2") on every data
# row; clone it from the row below onto the freshly inserted row.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)
$zj.Application.CutCopyMode = $false

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0.15

# ---------------------------------------------------------------------
# Keep the same tab active/selected as before the edit (the last quarter
# sheet, now shifted one slot to the right).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
